$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") cells hold numeric-looking text (t="inlineStr" in the
# source file). Excel auto-converts a numeric-looking string typed into a
# General-formatted cell into a real number, so we force these specific
# cells to Text format first (grouped into contiguous blocks so they all
# share one style) and then assign the new price strings. ---
$ws.Range("D2:D22").NumberFormat = "@"
$ws.Range("D24:D25").NumberFormat = "@"
$ws.Range("D40:D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49:D51").NumberFormat = "@"

# Rows 2-9: price-only updates
$ws.Range("D2").Value = "242.96"
$ws.Range("D3").Value = "23.18"
$ws.Range("D4").Value = "5.393"
$ws.Range("D5").Value = "0.05984"
$ws.Range("D6").Value = "3.403"
$ws.Range("D7").Value = "6.484"
$ws.Range("D8").Value = "0.8131"
$ws.Range("D9").Value = "0.8940"

# Rows 10-18: the ranking reshuffles - "One" jumps to rank 9 (row 10),
# pushing WazirX, MandalaExchangeToken, LiechtensteinCryptoassetsExchange,
# BitrueCoin, BitMartToken, MCDex, BitForexToken and CoinExToken each down
# one row, with new prices/volume labels throughout.
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.01122"
$ws.Range("E10").Value = "9OneONEBestin24h"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1412"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.07421"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "0.03363"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03073"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09330"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.865"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001584"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04646"
$ws.Range("E18").Value = "17CoinExTokenCET"

# Row 19 (TigerCash): price-only update
$ws.Range("D19").Value = "0.006078"

# Row 20 (HotbitToken): price update, and it's no longer "Best in 24h"
# (that label moved to row 10 / One)
$ws.Range("D20").Value = "0.005017"
$ws.Range("E20").Value = "19HotbitTokenHTB"

# Remaining scattered price-only updates
$ws.Range("D21").Value = "0.0009818"
$ws.Range("D22").Value = "0.00007800"
$ws.Range("D24").Value = "3.615"
$ws.Range("D25").Value = "2.160"
$ws.Range("D40").Value = "0.03887"
$ws.Range("D41").Value = "0.006233"
$ws.Range("D42").Value = "0.1073"
$ws.Range("D43").Value = "0.002800"
$ws.Range("D44").Value = "0.007189"
$ws.Range("D45").Value = "0.00005184"
$ws.Range("D47").Value = "0.0005798"
$ws.Range("D49").Value = "0.002299"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("D51").Value = "0.0002000"
